# Add a new "friendly" download-error message (Spanish/English) as a new
# row at the bottom of the Idiomas (languages) table, matching the style
# used by the other multi-line rows (e.g. the "Updated Version" row above).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$spanishText = "Algo ha ido mal. `nPor favor, inténtalo de nuevo más tarde. `nSi el problema persiste, contacta conmigo a través de mis redes sociales"
$englishText = "Something has gone wrong. `nPlease try again later. `nIf the problem persists, contact me through my social networks"

$newRow = $ws.Range("A24").Row + 1

$ws.Cells.Item($newRow, 1).Value = $spanishText
$ws.Cells.Item($newRow, 2).Value = $englishText

# Match formatting of the row above it (center aligned, word-wrapped, taller row).
$targetRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 2))
$targetRange.HorizontalAlignment = -4108
$targetRange.WrapText = $true
$ws.Rows.Item($newRow).RowHeight = 60

$ws.Cells.Item($newRow, 2).Select()
